$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.275.50"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3
$ws.Range("D3").Value = "1.594.27"
$ws.Range("E3").Value = "  +0.40%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
$ws.Range("E6").Value = "  +0.20%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.30%  "

# Row 9
$ws.Range("E9").Value = "  +0.58%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.24%  "

# Row 12
$ws.Range("D12").Value = "1.819.22"
$ws.Range("E12").Value = "  +0.42%  "

# Row 13
$ws.Range("D13").Value = "1.604.46"
$ws.Range("E13").Value = "  +1.08%  "

# Row 14
$ws.Range("E14").Value = "  -0.42%  "

# Row 15
$ws.Range("E15").Value = "  -1.41%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "

# Row 17
$ws.Range("D17").Value = "26.270.96"
$ws.Range("E17").Value = "  +0.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.63%  "

# Row 19
$ws.Range("E19").Value = "  +4.22%  "

# Row 20
$ws.Range("E20").Value = "  -0.39%  "

# Row 21
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "

# Row 23
$ws.Range("E23").Value = "  +2.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.41%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("E27").Value = "  +0.56%  "

# Row 28
$ws.Range("E28").Value = "  +0.31%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.40%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0492"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "

# Row 31
$ws.Range("E31").Value = "  +0.06%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.61%  "

# Row 33
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.469.54"
$ws.Range("E33").Value = "  +4.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.48%  "

# Row 35
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("E36").Value = "  +0.50%  "

# Row 37
$ws.Range("E37").Value = "  -3.12%  "

# Row 38
$ws.Range("E38").Value = "  -0.34%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.816"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "

# Row 40
$ws.Range("E40").Value = "  -1.88%  "

# Row 41
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("E42").Value = "  +1.19%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.931"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.64%  "

# Row 44
$ws.Range("D44").Value = "1.731.63"
$ws.Range("E44").Value = "  +0.45%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.755"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.23%  "

# Row 48
$ws.Range("E48").Value = "  -0.66%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0501"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "  +0.16%  "

# Row 51
$ws.Range("E51").Value = "  -0.11%  "
